{"js": "const body = context.document.body;\n\nconst results1 = body.search(\"wih to rip\", { matchCase: true, matchWholeWord: false });\nresults1.load(\"text\");\nawait context.sync();\n\nif (results1.items.length > 0) {\n  results1.items[0].insertText(\"wish to rip\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\nconst results2 = body.search(\"That is what this war is, Evil.\", { matchCase: true, matchWholeWord: false });\nresults2.load(\"text\");\nawait context.sync();\n\nif (results2.items.length > 0) {\n  results2.items[0].insertText(\"This is what this war is, Evil.\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $full = $p.Range.Text\n    $core = $full.Substring(0, $full.Length - 1)\n\n    if ($core.Contains(\"wih to rip\")) {\n        $p.Range.Text = $core.Replace(\"wih to rip\", \"wish to rip\")\n    }\n    elseif ($core -eq \"That is what this war is, Evil.\") {\n        $p.Range.Text = $core.Replace(\"That is what this war is, Evil.\", \"This is what this war is, Evil.\")\n    }\n}\n"}
